# Apply updated crypto price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''90.747.44'
$ws.Range("E2").Value = '''  +2.02%  '

# Row 3
$ws.Range("D3").Value = '''3.178.21'
$ws.Range("E3").Value = '''  -3.16%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '''  +0.07%  '

# Row 5
$ws.Range("D5").Value = '''212.28'
$ws.Range("E5").Value = '''  -0.86%  '

# Row 6
$ws.Range("D6").Value = '''614.98'
$ws.Range("E6").Value = '''  -2.44%  '

# Row 7
$ws.Range("D7").Value = '''0.387'
$ws.Range("E7").Value = '''  +0.96%  '

# Row 8
$ws.Range("E8").Value = '''  -0.55%  '

# Row 9
$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '''  +0.09%  '

# Row 10
$ws.Range("D10").Value = '''3.175.58'
$ws.Range("E10").Value = '''  -3.08%  '

# Row 11
$ws.Range("D11").Value = '''0.573'
$ws.Range("E11").Value = '''  -0.73%  '

# Row 12
$ws.Range("E12").Value = '''  -6.25%  '

# Row 13
$ws.Range("D13").Value = '''0.0000252'
$ws.Range("E13").Value = '''  -4.79%  '

# Row 14
$ws.Range("D14").Value = '''90.477.67'
$ws.Range("E14").Value = '''  +1.77%  '

# Row 15
$ws.Range("D15").Value = '''3.766.85'
$ws.Range("E15").Value = '''  -3.10%  '

# Row 16
$ws.Range("D16").Value = '''32.72'
$ws.Range("E16").Value = '''  -4.85%  '

# Row 17
$ws.Range("D17").Value = '''5.21'
$ws.Range("E17").Value = '''  -3.64%  '

# Row 18
$ws.Range("D18").Value = '''3.179.43'
$ws.Range("E18").Value = '''  -3.25%  '

# Row 19
$ws.Range("D19").Value = '''3.25'
$ws.Range("E19").Value = '''  +3.80%  '

# Row 20
$ws.Range("D20").Value = '''13.41'
$ws.Range("E20").Value = '''  -5.38%  '

# Row 21
$ws.Range("D21").Value = '''435.79'
$ws.Range("E21").Value = '''  -0.49%  '

# Row 22
$ws.Range("D22").Value = '''0.0000186'
$ws.Range("E22").Value = '''  +36.77%  '

# Row 23
$ws.Range("D23").Value = '''8.52'
$ws.Range("E23").Value = '''  -4.51%  '

# Row 24
$ws.Range("E24").Value = '''  -5.20%  '

# Row 25
$ws.Range("E25").Value = '''  -2.61%  '

# Row 26
$ws.Range("D26").Value = '''11.81'
$ws.Range("E26").Value = '''  -4.66%  '

# Row 27
$ws.Range("D27").Value = '''3.382.32'
$ws.Range("E27").Value = '''  -1.96%  '

# Row 28
$ws.Range("D28").Value = '''74.77'
$ws.Range("E28").Value = '''  -2.66%  '

# Row 29
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '''  +0.04%  '

# Row 30
$ws.Range("D30").Value = '''0.169'
$ws.Range("E30").Value = '''  -7.75%  '

# Row 31
$ws.Range("E31").Value = '''  +0.01%  '

# Row 32
$ws.Range("D32").Value = '''4.20'
$ws.Range("E32").Value = '''  +36.33%  '

# Row 33
$ws.Range("E33").Value = '''  -4.80%  '

# Row 34
$ws.Range("D34").Value = '''534.56'
$ws.Range("E34").Value = '''  -5.78%  '

# Row 35
$ws.Range("D35").Value = '''6.89'
$ws.Range("E35").Value = '''  -4.26%  '

# Row 36
$ws.Range("D36").Value = '''1.86'
$ws.Range("E36").Value = '''  -5.61%  '

# Row 37
$ws.Range("D37").Value = '''1.24'
$ws.Range("E37").Value = '''  -10.70%  '

# Row 38
$ws.Range("D38").Value = '''21.91'
$ws.Range("E38").Value = '''  -3.52%  '

# Row 39
$ws.Range("E39").Value = '''  +2.40%  '

# Row 40
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '''  +0.15%  '

# Row 41
$ws.Range("D41").Value = '''0.127'
$ws.Range("E41").Value = '''  -9.21%  '

# Row 42
$ws.Range("E42").Value = '''  -0.30%  '

# Row 43
$ws.Range("B43").Value = '''Stacks'
$ws.Range("C43").Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '''1.91'
$ws.Range("E43").Value = '''  -6.35%  '

# Row 44
$ws.Range("B44").Value = '''PolygonEcosystemToken'
$ws.Range("C44").Value = '''https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").Value = '''0.375'
$ws.Range("E44").Value = '''  -6.55%  '

# Row 45
$ws.Range("D45").Value = '''146.75'
$ws.Range("E45").Value = '''  -5.50%  '

# Row 46
$ws.Range("D46").Value = '''44.61'
$ws.Range("E46").Value = '''  -0.91%  '

# Row 47
$ws.Range("D47").Value = '''172.77'
$ws.Range("E47").Value = '''  -4.83%  '

# Row 48
$ws.Range("E48").Value = '''  -2.94%  '

# Row 49
$ws.Range("D49").Value = '''1.23'
$ws.Range("E49").Value = '''  -6.11%  '

# Row 50
$ws.Range("D50").Value = '''0.614'
$ws.Range("E50").Value = '''  -2.17%  '

# Row 51
$ws.Range("E51").Value = '''  -3.95%  '
